$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in column H, matching the style used by the other
# header cells (e.g. G1) so it reuses the same cell style index.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add the corresponding value for row 2 in the new Save column
$ws.Range("H2").Value = 0
